$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 453.41666
$ws.Range("I33").Value = 453.41666
$ws.Range("K33").Value = 453.41666
$ws.Range("M33").Value = -224.41666
$ws.Range("H40").Value = 1313.4286
$ws.Range("I40").Value = 1282.3334
$ws.Range("K40").Value = 1282.3334
$ws.Range("M40").Value = -1107.3334
$ws.Range("H46").Value = 4000
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H60").Value = 4000
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H101").Value = 263
$ws.Range("I101").Value = 269.2
$ws.Range("J101").Value = 232
$ws.Range("K101").Value = 807.5999999999999
$ws.Range("L101").Value = 696
$ws.Range("M101").Value = 814.4000000000001
$ws.Range("N101").Value = -3940
$ws.Range("H131").Value = 1244
$ws.Range("I131").Value = 1244
$ws.Range("K131").Value = 3732
$ws.Range("M131").Value = 1308

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1885.9333
$ws.Range("I74").Value = 2065.5
$ws.Range("K74").Value = 2065.5
$ws.Range("M74").Value = -1191.5
$ws.Range("H77").Value = 1885.9333
$ws.Range("I77").Value = 2065.5
$ws.Range("K77").Value = 10327.5
$ws.Range("M77").Value = -5959.5
$ws.Range("H132").Value = 2581.25
$ws.Range("I132").Value = 2620.1333
$ws.Range("J132").Value = 1998
$ws.Range("K132").Value = 7860.3999
$ws.Range("L132").Value = 5994
$ws.Range("M132").Value = -5330.3999
$ws.Range("N132").Value = -11054

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1979
$ws.Range("I86").Value = 2065.5
$ws.Range("J86").Value = 1633
$ws.Range("K86").Value = 2065.5
$ws.Range("L86").Value = 1633
$ws.Range("M86").Value = -942.5
$ws.Range("N86").Value = -3879
$ws.Range("H89").Value = 1979
$ws.Range("I89").Value = 2065.5
$ws.Range("J89").Value = 1633
$ws.Range("K89").Value = 10327.5
$ws.Range("L89").Value = 8165
$ws.Range("M89").Value = -4711.5
$ws.Range("N89").Value = -19397
$ws.Range("H107").Value = 1879.6
$ws.Range("I107").Value = 1849.5
$ws.Range("K107").Value = 1849.5
$ws.Range("M107").Value = 70.5
$ws.Range("H134").Value = 5031
$ws.Range("I134").Value = 5225
$ws.Range("K134").Value = 15675
$ws.Range("M134").Value = -13140

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 502282.94
$ws.Range("I19").Value = 588509
$ws.Range("J19").Value = 13668.667
$ws.Range("K19").Value = 588509
$ws.Range("L19").Value = 13668.667
$ws.Range("M19").Value = -588339
$ws.Range("N19").Value = -14008.667
$ws.Range("H24").Value = 502282.94
$ws.Range("I24").Value = 588509
$ws.Range("J24").Value = 13668.667
$ws.Range("K24").Value = 588509
$ws.Range("L24").Value = 13668.667
$ws.Range("M24").Value = -588339
$ws.Range("N24").Value = -14008.667
$ws.Range("H107").Value = 1319.875
$ws.Range("J107").Value = 3000
$ws.Range("L107").Value = 3000
$ws.Range("N107").Value = -6840

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 131
$ws.Range("I2").Value = 121
$ws.Range("J2").Value = 134.33333
$ws.Range("K2").Value = 726
$ws.Range("L2").Value = 805.9999799999999
$ws.Range("M2").Value = -613
$ws.Range("N2").Value = -1031.99998
$ws.Range("H12").Value = 363.54544
$ws.Range("J12").Value = 315
$ws.Range("L12").Value = 945
$ws.Range("N12").Value = -1291
$ws.Range("H34").Value = 49966.816
$ws.Range("J34").Value = 57748.156
$ws.Range("L34").Value = 173244.468
$ws.Range("N34").Value = -173412.468
$ws.Range("H35").Value = 2000
$ws.Range("I35").Value = 2000
$ws.Range("K35").Value = 6000
$ws.Range("M35").Value = -5712
$ws.Range("H38").Value = 242.33333
$ws.Range("I38").Value = 363.33334
$ws.Range("J38").Value = 121.333336
$ws.Range("K38").Value = 1090.00002
$ws.Range("L38").Value = 364.000008
$ws.Range("M38").Value = -743.0000199999999
$ws.Range("N38").Value = -1058.000008
$ws.Range("H39").Value = 8000.5
$ws.Range("J39").Value = 9000
$ws.Range("L39").Value = 27000
$ws.Range("N39").Value = -27588
$ws.Range("H40").Value = 20
$ws.Range("I40").Value = 20
$ws.Range("K40").Value = 80
$ws.Range("M40").Value = -11
$ws.Range("H44").Value = 999999
$ws.Range("I44").Value = 999999
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 2999997
$ws.Range("L44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -2999599
$ws.Range("H46").Value = 663.75
$ws.Range("I46").Value = 663.75
$ws.Range("K46").Value = 1991.25
$ws.Range("M46").Value = -1900.25
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H49").Value = 4500
$ws.Range("I49").Value = 4500
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 13500
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -13344
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H54").Value = 10004
$ws.Range("I54").Value = 10004
$ws.Range("K54").Value = 30012
$ws.Range("M54").Value = -29453
$ws.Range("H55").Value = 11399
$ws.Range("I55").Value = 600
$ws.Range("J55").Value = 14998.667
$ws.Range("K55").Value = 1800
$ws.Range("L55").Value = 44996.001
$ws.Range("M55").Value = -1623
$ws.Range("N55").Value = -45350.001
$ws.Range("H57").Value = 5555
$ws.Range("J57").Value = 5555
$ws.Range("L57").Value = 16665
$ws.Range("N57").Value = -17783
$ws.Range("H62").Value = 7570.3335
$ws.Range("J62").Value = 6055.5
$ws.Range("L62").Value = 18166.5
$ws.Range("N62").Value = -19538.5
$ws.Range("H64").Value = 353
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H65").Value = 7570.3335
$ws.Range("J65").Value = 6055.5
$ws.Range("L65").Value = 54499.5
$ws.Range("N65").Value = -61363.5
$ws.Range("H67").Value = 353
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H70").Value = 11995
$ws.Range("I70").Value = 9999.5
$ws.Range("J70").Value = 13325.333
$ws.Range("K70").Value = 29998.5
$ws.Range("L70").Value = 39975.999
$ws.Range("M70").Value = -29683.5
$ws.Range("N70").Value = -40605.999
$ws.Range("H73").Value = 11995
$ws.Range("I73").Value = 9999.5
$ws.Range("J73").Value = 13325.333
$ws.Range("K73").Value = 29998.5
$ws.Range("L73").Value = 39975.999
$ws.Range("M73").Value = -28906.5
$ws.Range("N73").Value = -42159.999
$ws.Range("H132").Value = 1557
$ws.Range("I132").Value = 1544.6364
$ws.Range("K132").Value = 13901.7276
$ws.Range("M132").Value = -11371.7276

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 15000000
$ws.Range("I18").Value = 15000000
$ws.Range("K18").Value = 15000000
$ws.Range("M18").Value = -14999707
$ws.Range("H27").Value = 40000
$ws.Range("J27").Value = 40000
$ws.Range("L27").Value = 40000
$ws.Range("N27").Value = -40332
$ws.Range("H70").Value = 33548.832
$ws.Range("I70").Value = 38758.7
$ws.Range("J70").Value = 7499.5
$ws.Range("K70").Value = 38758.7
$ws.Range("L70").Value = 7499.5
$ws.Range("M70").Value = -38488.7
$ws.Range("N70").Value = -8039.5
$ws.Range("H73").Value = 33548.832
$ws.Range("I73").Value = 38758.7
$ws.Range("J73").Value = 7499.5
$ws.Range("K73").Value = 38758.7
$ws.Range("L73").Value = 7499.5
$ws.Range("M73").Value = -37822.7
$ws.Range("N73").Value = -9371.5
$ws.Range("H122").Value = 1371.3158
$ws.Range("I122").Value = 1371.3158
$ws.Range("K122").Value = 4113.9474
$ws.Range("M122").Value = -1663.9474

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 3954.5
$ws.Range("J11").Value = 3954.5
$ws.Range("L11").Value = 3954.5
$ws.Range("N11").Value = -4234.5
$ws.Range("H25").Value = 8000
$ws.Range("I25").Value = 8000
$ws.Range("K25").Value = 8000
$ws.Range("M25").Value = -7770
$ws.Range("H122").Value = 5089.2
$ws.Range("I122").Value = 3841.8572
$ws.Range("K122").Value = 11525.5716
$ws.Range("M122").Value = -9075.571599999999
$ws.Range("H136").Value = 45456704
$ws.Range("I136").Value = 2064.5715
$ws.Range("J136").Value = 125002320
$ws.Range("K136").Value = 6193.7145
$ws.Range("L136").Value = 375006960
$ws.Range("M136").Value = -3643.7145
$ws.Range("N136").Value = -375012060

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5268.3887
$ws.Range("I136").Value = 1653.5
$ws.Range("K136").Value = 4960.5
$ws.Range("M136").Value = -2410.5
